# Commit: the deck's two theme parts are swapped -
#   ppt/theme/theme1.xml  (the slide master's theme, "Integral" / Red
#                           Violet colour scheme)
#   ppt/theme/theme2.xml  (the notes master's theme, "Office Theme" /
#                           Office colour scheme)
# end up carrying each other's colours - i.e. theme1.xml becomes the
# "Office" palette (theme2.xml's old content).
#
# PowerPoint's object model exposes the slide master's theme colours
# through SlideMaster.Theme.ThemeColorScheme: a 12-slot collection
# (msoThemeColorSchemeIndex order - dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink) whose .RGB is a normal OLE_COLOR (0x00BBGGRR).
# Re-pointing every slot to the "Office" palette reproduces the new
# theme1.xml colour scheme.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$theme.Name = "Office Theme"

$colors = $theme.ThemeColorScheme
$colors.Name = "Office"

function Set-ThemeRGB($index, $r, $g, $b) {
    $rgb = $r + ($g * 256) + ($b * 65536)
    $colors.Item($index).RGB = $rgb
}

# msoThemeColorSchemeIndex slots: 1=dk1 2=lt1 3=dk2 4=lt2
# 5-10=accent1..accent6 11=hlink 12=folHlink
Set-ThemeRGB 1  0x00 0x00 0x00   # dk1
Set-ThemeRGB 2  0xFF 0xFF 0xFF   # lt1
Set-ThemeRGB 3  0x44 0x54 0x6A   # dk2
Set-ThemeRGB 4  0xE7 0xE6 0xE6   # lt2
Set-ThemeRGB 5  0x5B 0x9B 0xD5   # accent1
Set-ThemeRGB 6  0xED 0x7D 0x31   # accent2
Set-ThemeRGB 7  0xA5 0xA5 0xA5   # accent3
Set-ThemeRGB 8  0xFF 0xC0 0x00   # accent4
Set-ThemeRGB 9  0x44 0x72 0xC4   # accent5
Set-ThemeRGB 10 0x70 0xAD 0x47   # accent6
Set-ThemeRGB 11 0x05 0x63 0xC1   # hlink
Set-ThemeRGB 12 0x95 0x4F 0x72   # folHlink
